$d = $word.ActiveDocument

# Paragraphs in the ingredients list originally begin with two tabs and a
# space (produced by two separate runs: one run with a single <w:tab/>,
# and a second run with <w:tab/> followed by a literal space). The edit
# collapses that leading "\t\t " prefix into a single run containing 16
# literal space characters, leaving the remainder of the paragraph (the
# ingredient text, in its own run) untouched.
$newPrefix = "                "  # 16 spaces

# Work from the end towards the start so that replacing a paragraph's
# content does not disturb the indices/positions of paragraphs not yet
# processed.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $range = $p.Range
    $text = $range.Text

    if ($text.StartsWith("`t`t ")) {
        # Strip the leading two tabs + space, and drop the trailing
        # paragraph-mark character(s) captured by Range.Text.
        $rest = $text.Substring(3)
        $rest = $rest.TrimEnd([char]13, [char]7)

        # Escape for safe embedding inside the XML we hand to InsertXML.
        $restEscaped = $rest.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

        $newParagraphXml = '<w:p><w:r><w:t xml:space="preserve">' + $newPrefix + '</w:t></w:r><w:r><w:t>' + $restEscaped + '</w:t></w:r></w:p>'

        $range.InsertXML($newParagraphXml)
    }
}
